# Tabela de precos atualizada em 16/12/2023
# The exam "PESQUISA DE FUNGOS" (row 296) was removed from the price list.
# Deleting the row shifts every following row up by one and drops the
# now-unused shared string; reapply the existing value-sort over the
# shrunk range so the stored sort state stays consistent with the data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "PESQUISA DE FUNGOS" row entirely (cells shift up).
$ws.Rows.Item(296).Delete()

# Re-apply the ascending sort on column B over the new (one row shorter)
# range so the worksheet's remembered sort range matches the data again.
$sortObj = $ws.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws.Range("B2:B720"), 0, 1, 0, 0)
$sortObj.SetRange($ws.Range("A2:C720"))
$sortObj.Header = 2
$sortObj.MatchCase = $false
$sortObj.Orientation = 1
$sortObj.Apply()
